$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 152; existing rows 152..242 shift down to 153..243
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with its data (columns A-T)
$ws.Range("A152").Value = 10
$ws.Range("B152").Value = "Vega Modelo de Temuco"
$ws.Range("C152").Value = "La Araucanía"
$ws.Range("D152").Value = 44606
$ws.Range("D152").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E152").Value = 9
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100101
$ws.Range("H152").Value = "Berries"
$ws.Range("I152").Value = 100112025
$ws.Range("J152").Value = "Frutilla"
$ws.Range("K152").Value = "Sin especificar"
$ws.Range("L152").Value = "Primera"
$ws.Range("M152").Value = 50
$ws.Range("N152").Value = 8000
$ws.Range("O152").Value = 8000
$ws.Range("P152").Value = 8000
$ws.Range("Q152").Value = "$/caja 7 kilos"
$ws.Range("R152").Value = "Región de La Araucanía"
$ws.Range("S152").Value = 1143
$ws.Range("T152").Value = 7
